$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# This team-stats export has a "Date" column whose values were written as
# "11-22-2012-13" (day-month concatenated with the season label) instead of
# the actual game date. Because the underlying NBA stats page showed the
# date one day off, every row in the column needs to be corrected to the
# real game date "2012-11-22".
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# Locate the "Date" header so the fix targets the right column even if the
# sheet layout shifts.
$dateCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Text -eq "Date") {
        $dateCol = $c
        break
    }
}
if ($dateCol -eq 0) {
    $dateCol = 58  # fall back to column BF
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Text -eq "11-22-2012-13") {
        $origStyle = $cell.Style
        # Force the cell to Text format first so Excel stores the literal
        # string "2012-11-22" instead of re-interpreting it as a date
        # serial number, then restore the original (default) cell style so
        # no stray formatting is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = "2012-11-22"
        $cell.Style = $origStyle
    }
}
